$wb = $excel.ActiveWorkbook

# --- Sheet2: insert a new header row at the top, fill it in ---
$ws2 = $wb.Worksheets.Item("Sheet2")

# Drop the existing hyperlinks now (their cell refs are about to shift down a
# row); we'll re-add them at the new locations right after the insert.
$ws2.Range("A1").Hyperlinks.Delete() | Out-Null

$ws2.Rows.Item(1).Insert()

$ws2.Cells.Item(1,1).Value = "Name"
$ws2.Cells.Item(1,2).Value = "Lastname"
$ws2.Cells.Item(1,3).Value = "email"
$ws2.Cells.Item(1,4).Value = "telefone"
$ws2.Cells.Item(1,5).Value = "password"
$ws2.Cells.Item(1,6).Value = "confirm password"

# Re-create the hyperlinks shifted one row down, matching the original
# targets. (Hyperlinks.Add always re-applies its own "Hyperlink" style xf, so
# re-apply the built-in style afterwards to match the existing C2:C5 cells.)
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,3), "mailto:you@filthy.com") | Out-Null
$ws2.Cells.Item(2,3).Style = "Hyperlink"

$ws2.Hyperlinks.Add($ws2.Cells.Item(3,3), "mailto:doofenshmirtz@evilInc.com") | Out-Null
$ws2.Cells.Item(3,3).Style = "Hyperlink"

$ws2.Hyperlinks.Add($ws2.Cells.Item(4,3), "mailto:bigbills@carsale.com") | Out-Null
$ws2.Cells.Item(4,3).Style = "Hyperlink"

# D4 keeps its own visible text ("no cabe el link al video :c"); only the
# hyperlink's display/tooltip should be the YouTube URL, so restore the cell
# text after Add (which otherwise overwrites it with TextToDisplay).
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,4), "https://youtu.be/SBs455jwb8w?si=EOOIHk3wa9iFf3z5", "", "", "https://youtu.be/SBs455jwb8w?si=EOOIHk3wa9iFf3z5") | Out-Null
$ws2.Cells.Item(4,4).Value = "no cabe el link al video :c"
$ws2.Cells.Item(4,4).Style = "Hyperlink"

$ws2.Hyperlinks.Add($ws2.Cells.Item(5,3), "mailto:engineering@mylimit.com") | Out-Null
$ws2.Cells.Item(5,3).Style = "Hyperlink"

# New "confirm password" column needs its own width (stored width 15.5).
$ws2.Columns.Item(6).ColumnWidth = 14.6

# --- Sheet1: move the selection ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B10").Select() | Out-Null

# --- Sheet2 becomes the active sheet, with a new selection ---
$ws2.Activate()
$ws2.Range("D14").Select() | Out-Null
